$wb = $excel.ActiveWorkbook

# Sheet 1: illnessratio - new values, and row 7 removed (6 rows only)
$ws1 = $wb.Worksheets.Item("illnessratio")
$ws1.Range("A1").Value = 4.85710030105084
$ws1.Range("A2").Value = 2.29806616604013
$ws1.Range("A3").Value = 7.23794521699157
$ws1.Range("A4").Value = 3.93434262197864
$ws1.Range("A5").Value = 4.76380389296696
$ws1.Range("A6").Value = 1.99203903661893
$ws1.Range("A7").ClearContents()

# Sheet 2: illnessday - new values, still 7 rows
$ws2 = $wb.Worksheets.Item("illnessday")
$ws2.Range("A1").Value = 4.87477027960459
$ws2.Range("A2").Value = 2.30412078878102
$ws2.Range("A3").Value = 7.38847255076552
$ws2.Range("A4").Value = 4.30866246722926
$ws2.Range("A5").Value = 1.88120132825269
$ws2.Range("A6").Value = 3.56397999754844
$ws2.Range("A7").Value = 2.09056790748145

# Sheet 3: chronicratio - new values, still 7 rows
$ws3 = $wb.Worksheets.Item("chronicratio")
$ws3.Range("A1").Value = 4.87477027960459
$ws3.Range("A2").Value = 2.30412078878102
$ws3.Range("A3").Value = 7.3884725507655
$ws3.Range("A4").Value = 4.30866246722926
$ws3.Range("A5").Value = 1.88120132825269
$ws3.Range("A6").Value = 3.56397999754844
$ws3.Range("A7").Value = 2.09056790748145
